# Case 380 kV: update pl_mw.xlsx line-flow results (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:D
$rowsBlock1 = @(
  @(1.357163811744442, 0.06271765734092583, 0.01862097285227904),
  @(1.327593896070084, 0.05710098048051293, 0.01843714672380159),
  @(1.310144512444339, 0.05368371984020825, 0.01833035663716842),
  @(1.30321184614121, 0.05229899887002887, 0.01828838046023584),
  @(1.302071450880732, 0.05206953893669208, 0.01828150383940752),
  @(1.310050294430567, 0.0536650133364418, 0.01832978427205845),
  @(1.346821664682295, 0.06077447438408967, 0.01855633526195533),
  @(1.424525056529745, 0.07496852607511073, 0.01904832260094835),
  @(1.48501619201781, 0.08555637250196924, 0.01943824263862837),
  @(1.513272789650557, 0.09040895694988649, 0.01962167843936413),
  @(1.524078808906665, 0.09225177728504264, 0.01969200061473231),
  @(1.521746837395597, 0.09185465844817031, 0.01967681746494776),
  @(1.514159688013763, 0.09056046129950346, 0.01962744674099071),
  @(1.50952611436594, 0.08976841385410239, 0.01959731725527192),
  @(1.48318439154022, 0.0852399764660845, 0.01942637554961024),
  @(1.467213562989343, 0.08247121912985733, 0.01932305145854585),
  @(1.458097134548893, 0.08088209945189817, 0.01926419321913997),
  @(1.455022429141678, 0.08034463239305012, 0.01924436328646095),
  @(1.468906487585286, 0.082765606152023, 0.01933399147274173),
  @(1.516385347818414, 0.09094045537622719, 0.01964192490303418),
  @(1.548032375669891, 0.09631383117434211, 0.01984817864663668),
  @(1.531085447948101, 0.09344313668975701, 0.01973764352162632),
  @(1.46814091258679, 0.08263250545049061, 0.01932904379667733),
  @(1.402906134713305, 0.07110105832022384, 0.01891018445476789)
)
$arrBlock1 = New-Object 'object[,]' 24,3
for ($r = 0; $r -lt 24; $r++) {
  for ($c = 0; $c -lt 3; $c++) {
    $arrBlock1[$r,$c] = $rowsBlock1[$r][$c]
  }
}
$ws.Range("B2:D25").Value = $arrBlock1

# Columns F:G
$rowsBlock2 = @(
  @(3.640032366702087, 0.002620915803398652),
  @(3.640752705795734, 0.002624731211801146),
  @(3.642527755172551, 0.002627199957784457),
  @(3.643586497747066, 0.002628237791056885),
  @(3.643782566032684, 0.002628412046115081),
  @(3.642540675400696, 0.002627213825292146),
  @(3.640004275178782, 0.002622205249689759),
  @(3.645597043255094, 0.002613379267756412),
  @(3.656141697033434, 0.002607495615162639),
  @(3.662335308522643, 0.002604948109315341),
  @(3.664881374049131, 0.002604001883463658),
  @(3.66432411080811, 0.002604204850706073),
  @(3.662540753960457, 0.002604869893388464),
  @(3.661474525545827, 0.002605279652506215),
  @(3.655765021720327, 0.002607664688596793),
  @(3.652620050754322, 0.002609160805612576),
  @(3.650942629998156, 0.00261003347997569),
  @(3.650397272060516, 0.002610331041506572),
  @(3.65294123154284, 0.002609000284889245),
  @(3.663059124116558, 0.002604674054083268),
  @(3.670841384570508, 0.002601954165320701),
  @(3.666580882009924, 0.002603396009756431),
  @(3.652795618665564, 0.002609072817414409),
  @(3.642953462876747, 0.002615660972355677)
)
$arrBlock2 = New-Object 'object[,]' 24,2
for ($r = 0; $r -lt 24; $r++) {
  for ($c = 0; $c -lt 2; $c++) {
    $arrBlock2[$r,$c] = $rowsBlock2[$r][$c]
  }
}
$ws.Range("F2:G25").Value = $arrBlock2

# Columns I:N
$rowsBlock3 = @(
  @(2.296344833468481, 0.1643149148961385, 1.159734745125633, 0.368301654903604, 0.3640308269039387, 3.887144544557628),
  @(2.300410310924377, 0.1648918657895084, 1.125907901754431, 0.3670671761061541, 0.3590985506017503, 3.90651228458281),
  @(2.303679067529821, 0.1652705442291964, 1.105778963716205, 0.3664500298030759, 0.3562439196994447, 3.919274239014364),
  @(2.305205510935991, 0.1654310149011113, 1.097737663646711, 0.3662340324345337, 0.3551244396639035, 3.924693693281242),
  @(2.305470721000908, 0.1654580331452546, 1.096412167891799, 0.3662003124613875, 0.3549411995494474, 3.925606809746405),
  @(2.303698866390121, 0.1652726834488281, 1.105669862072062, 0.3664469729605599, 0.356228644501293, 3.919346441463858),
  @(2.297586349157903, 0.1645087855664915, 1.14793837088979, 0.3678468256948548, 0.3622941619042521, 3.893641994316283),
  @(2.291724889108387, 0.1632040088226248, 1.23590674676035, 0.37170623343534, 0.3755643528385377, 3.850136450756423),
  @(2.2911486133886, 0.1623623567379671, 1.303635797544416, 0.3752176072948004, 0.3861496974176575, 3.822375788975847),
  @(2.291695786836868, 0.1620046910724717, 1.335121333817256, 0.3769611073328036, 0.3911462028958681, 3.810658588247961),
  @(2.292019284495737, 0.1618728635652023, 1.34714108463865, 0.3776422733578784, 0.393064228663178, 3.806352594869438),
  @(2.291944442559895, 0.161901094488984, 1.344548112067258, 0.3774946419576679, 0.3926499945433193, 3.807274138498087),
  @(2.291720070937373, 0.161993773207568, 1.336108266059313, 0.3770167280342633, 0.3913034803474886, 3.810301704904248),
  @(2.291597779399467, 0.1620510117297442, 1.330951223878458, 0.376726716738645, 0.3904820791359569, 3.812173244730189),
  @(2.291129132672637, 0.1623862372228722, 1.301591711014197, 0.3751065999536678, 0.3858268016199347, 3.823159873966318),
  @(2.291048859214833, 0.1625983346002897, 1.283753378171554, 0.3741500905094739, 0.3830172756188546, 3.830133222881415),
  @(2.291078870351889, 0.1627227006231653, 1.273556812885914, 0.3736136932738248, 0.3814183696909126, 3.834229879815439),
  @(2.291102118629638, 0.1627652167883014, 1.270115356947002, 0.3734344443264774, 0.3808799404661372, 3.835631668273592),
  @(2.291049520494639, 0.1625755109397158, 1.28564572143398, 0.3742504888483609, 0.3833145895990384, 3.829382020871847),
  @(2.291782818773385, 0.1619664532675316, 1.338584626511476, 0.3771565352333539, 0.3916982804852438, 3.809408878270787),
  @(2.292939890019724, 0.1615894511107463, 1.373747735691609, 0.3791778258437688, 0.3973287638286251, 3.797119182153637),
  @(2.292260346267483, 0.1617887418771407, 1.35492895129525, 0.3780878853546596, 0.3943098603628243, 3.803608519828828),
  @(2.291048984273552, 0.1625858219515184, 1.284790009465695, 0.3742050566634276, 0.3831801230810896, 3.829721366452645),
  @(2.292655244046443, 0.1635363853849068, 1.211564949657799, 0.3705432254981318, 0.3718274405403434, 3.861167509497591)
)
$arrBlock3 = New-Object 'object[,]' 24,6
for ($r = 0; $r -lt 24; $r++) {
  for ($c = 0; $c -lt 6; $c++) {
    $arrBlock3[$r,$c] = $rowsBlock3[$r][$c]
  }
}
$ws.Range("I2:N25").Value = $arrBlock3

Write-Host "Case 380 kV pl_mw.xlsx results updated"